$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data values F2:F90
$timeTaken = @(
    "2021-10-05 13:40:03.931921",
    "2021-10-05 13:40:03.931935",
    "2021-10-05 13:40:03.931938",
    "2021-10-05 13:40:03.931941",
    "2021-10-05 13:40:03.931945",
    "2021-10-05 13:40:03.931948",
    "2021-10-05 13:40:03.931951",
    "2021-10-05 13:40:03.931954",
    "2021-10-05 13:40:03.931958",
    "2021-10-05 13:40:03.931961",
    "2021-10-05 13:40:03.931964",
    "2021-10-05 13:40:03.931966",
    "2021-10-05 13:40:03.931969",
    "2021-10-05 13:40:03.931972",
    "2021-10-05 13:40:03.931975",
    "2021-10-05 13:40:03.931978",
    "2021-10-05 13:40:03.931981",
    "2021-10-05 13:40:03.931985",
    "2021-10-05 13:40:03.931988",
    "2021-10-05 13:40:03.931991",
    "2021-10-05 13:40:03.931993",
    "2021-10-05 13:40:03.931996",
    "2021-10-05 13:40:03.931999",
    "2021-10-05 13:40:03.932002",
    "2021-10-05 13:40:03.932006",
    "2021-10-05 13:40:03.932009",
    "2021-10-05 13:40:03.932012",
    "2021-10-05 13:40:03.932015",
    "2021-10-05 13:40:03.932018",
    "2021-10-05 13:40:03.932021",
    "2021-10-05 13:40:03.932024",
    "2021-10-05 13:40:03.932027",
    "2021-10-05 13:40:03.932030",
    "2021-10-05 13:40:03.932033",
    "2021-10-05 13:40:03.932036",
    "2021-10-05 13:40:03.932039",
    "2021-10-05 13:40:03.932042",
    "2021-10-05 13:40:03.932045",
    "2021-10-05 13:40:03.932048",
    "2021-10-05 13:40:03.932051",
    "2021-10-05 13:40:03.932055",
    "2021-10-05 13:40:03.932058",
    "2021-10-05 13:40:03.932061",
    "2021-10-05 13:40:03.932064",
    "2021-10-05 13:40:03.932067",
    "2021-10-05 13:40:03.932070",
    "2021-10-05 13:40:03.932072",
    "2021-10-05 13:40:03.932075",
    "2021-10-05 13:40:03.932078",
    "2021-10-05 13:40:03.932081",
    "2021-10-05 13:40:03.932084",
    "2021-10-05 13:40:03.932087",
    "2021-10-05 13:40:03.932091",
    "2021-10-05 13:40:03.932094",
    "2021-10-05 13:40:03.932097",
    "2021-10-05 13:40:03.932100",
    "2021-10-05 13:40:03.932103",
    "2021-10-05 13:40:03.932106",
    "2021-10-05 13:40:03.932109",
    "2021-10-05 13:40:03.932112",
    "2021-10-05 13:40:03.932115",
    "2021-10-05 13:40:03.932118",
    "2021-10-05 13:40:03.932121",
    "2021-10-05 13:40:03.932124",
    "2021-10-05 13:40:03.932128",
    "2021-10-05 13:40:03.932132",
    "2021-10-05 13:40:03.932135",
    "2021-10-05 13:40:03.932138",
    "2021-10-05 13:40:03.932141",
    "2021-10-05 13:40:03.932144",
    "2021-10-05 13:40:03.932147",
    "2021-10-05 13:40:03.932150",
    "2021-10-05 13:40:03.932153",
    "2021-10-05 13:40:03.932156",
    "2021-10-05 13:40:03.932159",
    "2021-10-05 13:40:03.932162",
    "2021-10-05 13:40:03.932167",
    "2021-10-05 13:40:03.932170",
    "2021-10-05 13:40:03.932173",
    "2021-10-05 13:40:03.932177",
    "2021-10-05 13:40:03.932180",
    "2021-10-05 13:40:03.932183",
    "2021-10-05 13:40:03.932186",
    "2021-10-05 13:40:03.932189",
    "2021-10-05 13:40:03.932192",
    "2021-10-05 13:40:03.932195",
    "2021-10-05 13:40:03.932198",
    "2021-10-05 13:40:03.932201",
    "2021-10-05 13:40:03.932204"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

